# Add 4 new quote rows (episode 17) to the bottom of the "Worksheet" sheet,
# matching: text (A) / character (B) / season (C) / episode (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A237").Value = "Seeing that you recognize this sword, you must be a mage."
$ws.Range("B237").Value = "Jang Uk"
$ws.Range("C237").Value = 1
$ws.Range("D237").Value = 17

$ws.Range("A238").Value = "I did not expect to run into anyone who recognized me at a place like this."
$ws.Range("B238").Value = "Jang Uk"
$ws.Range("C238").Value = 1
$ws.Range("D238").Value = 17

$ws.Range("A239").Value = "Are you saying the great hero who destroyed the ice stone to save the world was tempted to use it just to save someone's life?"
$ws.Range("B239").Value = "Prince Go Won"
$ws.Range("C239").Value = 1
$ws.Range("D239").Value = 17

$ws.Range("A240").Value = "It shows the foolishness and heartache of a man who was willing to go against the laws of nature to save a loved one."
$ws.Range("B240").Value = "Jang Uk"
$ws.Range("C240").Value = 1
$ws.Range("D240").Value = 17

# Match the author's final cursor position/selection on save.
$ws.Range("A239").Select()
